$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (product name casing) corrections - rows 11 and 12 swap
$ws.Range("C11").Value = 'SIG-3W Lilliput LED Torch &amp; Table Lamp'
$ws.Range("C12").Value = 'SIG-3w Lilliput LED Torch &amp; Table Lamp'

# Numeric value corrections (quantities, totals, subtotals, grand totals)
$numericUpdates = @{
    "B11" = 59408
    "D11" = 388.17
    "E11" = 463.78
    "F11" = 26
    "G11" = 10092.42
    "B12" = 47438
    "D12" = 401.81
    "E12" = 480.05
    "F12" = 2
    "G12" = 803.62
    "F51" = 73
    "G51" = 1919.17
    "B55" = 4494.85
    "F130" = 91
    "G130" = 4502.68
    "F133" = 13
    "G133" = 550.16
    "B140" = 58746.77
    "F183" = 2
    "G183" = 2045.22
    "B188" = 207596.99
    "F193" = 0
    "G193" = 0
    "B205" = 7329.96
    "F300" = 39
    "G300" = 11534.25
    "F305" = 30
    "G305" = 3111
    "F307" = 164
    "G307" = 3314.44
    "F308" = 24
    "G308" = 1968.72
    "F325" = 44
    "G325" = 4881.36
    "B331" = 213068.6
    "F341" = 235
    "G341" = 40605.65
    "F342" = 4
    "G342" = 691.16
    "F350" = 243
    "G350" = 17911.53
    "F352" = 97
    "G352" = 6738.59
    "F360" = 98
    "G360" = 14066.92
    "F361" = 8
    "G361" = 902.72
    "F362" = 60
    "G362" = 3958.8
    "F364" = 27
    "G364" = 3873.96
    "F365" = 41
    "G365" = 3406.28
    "F366" = 57
    "G366" = 2715.48
    "F369" = 9
    "G369" = 122.04
    "F400" = 25
    "G400" = 3866.5
    "F415" = 49
    "G415" = 4071.9
    "F417" = 645
    "G417" = 110507.85
    "F418" = 233
    "G418" = 35222.61
    "F419" = 4
    "G419" = 1673.96
    "F421" = 17
    "G421" = 2725.44
    "F426" = 24
    "G426" = 852.96
    "F428" = 80
    "G428" = 1625.6
    "F429" = 444
    "G429" = 26404.68
    "F430" = 441
    "G430" = 18160.38
    "F434" = 136
    "G434" = 19537.76
    "B435" = 691283.24
    "F437" = 123
    "G437" = 22580.34
    "B453" = 109982.5
    "F465" = 42
    "G465" = 3093.3
    "B477" = 34787.5
    "F482" = 197
    "G482" = 2094.11
    "F485" = 79
    "G485" = 3313.26
    "F490" = 137
    "G490" = 5693.72
    "F494" = 92
    "G494" = 604.4400000000001
    "F497" = 59
    "G497" = 8734.360000000001
    "F498" = 18
    "G498" = 1193.22
    "B507" = 123992.58
    "F558" = 9
    "G558" = 1658.07
    "F560" = 51
    "G560" = 2226.15
    "F563" = 33
    "G563" = 3039.96
    "F564" = 47
    "G564" = 4834.89
    "F568" = 44
    "G568" = 3166.68
    "F569" = 7
    "G569" = 1121.33
    "F576" = 9
    "G576" = 323.01
    "B577" = 37558.45
    "F702" = 0
    "G702" = 0
    "B704" = 19358.04
    "F799" = 117
    "G799" = 10734.75
    "B807" = 57421.98
    "F843" = 269
    "G843" = 21939.64
    "F846" = 88
    "G846" = 13597.76
    "F847" = 186
    "G847" = 15170.16
    "F848" = 361
    "G848" = 48049.1
    "F852" = 125
    "G852" = 2715
    "F853" = 186
    "G853" = 6934.08
    "F863" = 427
    "G863" = 61488
    "F865" = 298
    "G865" = 35971.58
    "B867" = 469723.73
    "F877" = 31
    "G877" = 1063.61
    "B878" = 3873.52
    "F893" = 8
    "G893" = 42173.04
    "F899" = 7
    "G899" = 99093.12
    "B902" = 433569.15
    "F909" = 54
    "G909" = 2466.18
    "F912" = 1784
    "G912" = 290988.24
    "B918" = 324502.53
    "B930" = 5942329.41
    "B931" = 5942329.41
}

foreach ($cellRef in $numericUpdates.Keys) {
    $ws.Range($cellRef).Value = $numericUpdates[$cellRef]
}

Write-Output "Applied $($numericUpdates.Count) numeric updates and 2 text updates."